# Section 9 (Refining Selections / LIKE) edit.
#
# Part 1: several existing bullet paragraphs get their single <w:r> run
# split into multiple runs interleaved with <w:proofErr> spell/grammar
# markers (purely cosmetic - visible text is unchanged). We rebuild each
# such paragraph's run content in place via Range.InsertXML so the exact
# w:proofErr markup from the target OOXML is reproduced.
#
# Part 2: a new "LIKE" bullet block (11 paragraphs) is appended at the
# end of the document.

$d = $word.ActiveDocument

function Set-ParaRuns($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $start = $p.Range.Start
    $end = $p.Range.End
    $r = $d.Range($start, $end)
    $r.Text = ""
    $r2 = $d.Range($start, $start)
    $full = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:p>'
    $r2.InsertXML($full)
}

# Each call below rewrites ONE paragraph's runs (paragraph mark / pPr /
# numbering stay untouched because only the paragraph's interior range is
# cleared before the replacement fragment is inserted at that same spot).
# Paragraph indices refer to positions in the ORIGINAL document - none of
# these calls change the total paragraph count, so the numbering stays
# valid across the whole sequence.

Set-ParaRuns 2 '<w:r><w:t xml:space="preserve">Distinct – </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>doesn’t</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> list any duplicates</w:t></w:r>'
Set-ParaRuns 3 '<w:r><w:t xml:space="preserve">SELECT DISTINCT </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>author_lname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> FROM </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>books;</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Set-ParaRuns 5 '<w:r><w:t xml:space="preserve">SELECT DISTINCT </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>author_fname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>author_lname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> FROM </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>books;</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Set-ParaRuns 8 '<w:r><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>author_lname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> FROM books ORDER BY </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>author_</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>lname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>;</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Set-ParaRuns 9 '<w:r><w:t>“</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>DESC”  =</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> descending</w:t></w:r>'
Set-ParaRuns 11 '<w:r><w:t xml:space="preserve">SELECT title, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>author_fname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>author_lname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> FROM books ORDER BY </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>2;</w:t></w:r><w:proofErr w:type="gramEnd"/>'
Set-ParaRuns 13 '<w:r><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>author_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>author_lname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> FROM books ORDER BY </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>author_lname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>author_fname</w:t></w:r><w:proofErr w:type="spellEnd"/>'
Set-ParaRuns 16 '<w:r><w:t xml:space="preserve">SELECT title FROM books ORDER BY </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stock_quantity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> LIMIT 5</w:t></w:r>'
Set-ParaRuns 19 '<w:r><w:t xml:space="preserve">To select from a starting point to the end, you </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>have to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> just use a gigantic number as the limit</w:t></w:r>'
Set-ParaRuns 20 '<w:r><w:t xml:space="preserve">IE: SELECT title FROM books LIMIT 5, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>18446744073709551615;</w:t></w:r><w:proofErr w:type="gramEnd"/>'

# Append new paragraphs (Section 9 - LIKE) at the end of the document
$endRange = $d.Content
$endRange.Collapse(0)
$newParasXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>LIKE – allows us to perform better searching for our data</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">WHERE allowed us to find exact matches, LIKE is like ‘contains’ or ‘starts with’. Use in combination. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">SELECT title </w:t></w:r><w:r><w:t>FROM books</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">WHERE </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>author_fname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> LIKE ‘%da%</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>’</w:t></w:r><w:r><w:t>;</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>This outputs any books with author first name that contains ‘da’</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>‘</w:t></w:r><w:r><w:t>%</w:t></w:r><w:r><w:t>’</w:t></w:r><w:r><w:t xml:space="preserve"> symbols are known as ‘wild cards’ </w:t></w:r><w:r><w:t>for any amount of characters</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Kinda</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> like regular expressions</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">‘_’ symbols are for ‘wild cards’ for a specific amount of characters. Each ‘_’ = 1 wild card character. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>It’s</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> also a way to match specific patterns. </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Ex: SELECT number FROM </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>phone_book</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> WHERE number LIKE (__</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>_)_</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>__-____</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Use an escape character to include actual percentage signs or underscores in searches</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>SELECT title FROM books WHERE title LIKE ‘%\%%’</w:t></w:r></w:p>'
$endRange.InsertXML($newParasXml)

